$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price/Volume columns so numeric-looking strings
# (e.g. "47.80") are preserved exactly as text, matching source data.
$textCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","E8","E9","D10","E10","D11","E11","D12","E12","D14","E14","D15","E15","E16","D17","E17","D18","E18","D19","E19","E20","E21","E22","D23","E23","D24","E24","E25","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","E33","E34","D35","E35","D36","E36","D37","E37","E38","E39","D40","E40","D42","E42","D43","E43","E44","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns
$ws.Range("D2").Value = '39.985.20'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '2.215.37'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '290.03'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '88.18'
$ws.Range("E6").Value = '  +5.53%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +1.41%  '
$ws.Range("D10").Value = '30.72'
$ws.Range("E10").Value = '  +3.80%  '
$ws.Range("D11").Value = '0.0781'
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '47.80'
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D14").Value = '6.48'
$ws.Range("E14").Value = '  +3.25%  '
$ws.Range("D15").Value = '2.557.93'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '2.199.90'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = '0.729'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").Value = '39.937.31'
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  +14.24%  '
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").Value = '65.63'
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").Value = '234.72'
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  +2.11%  '
$ws.Range("D27").Value = '1.83'
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  +6.34%  '
$ws.Range("D29").Value = '22.63'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").Value = '9.23'
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("D31").Value = '153.55'
$ws.Range("E31").Value = '  +2.93%  '
$ws.Range("D32").Value = '32.20'
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  +2.74%  '
$ws.Range("D35").Value = '0.0721'
$ws.Range("E35").Value = '  +2.57%  '
$ws.Range("D36").Value = '2.39'
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = '2.83'
$ws.Range("E37").Value = '  +6.55%  '
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("E39").Value = '  +3.12%  '
$ws.Range("D40").Value = '15.92'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D42").Value = '3.85'
$ws.Range("E42").Value = '  +5.39%  '
$ws.Range("D43").Value = '2.100.91'
$ws.Range("E43").Value = '  +8.75%  '
$ws.Range("E44").Value = '  +1.77%  '
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").Value = '9.97'
$ws.Range("E46").Value = '  +6.57%  '
$ws.Range("D47").Value = '17.82'
$ws.Range("E47").Value = '  +9.95%  '
$ws.Range("D48").Value = '2.67'
$ws.Range("E48").Value = '  +2.79%  '
$ws.Range("D49").Value = '2.431.08'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '69.46'
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("D51").Value = '1.45'
$ws.Range("E51").Value = '  +3.56%  '

# Update Coin (B) and Link (C) columns (row reorderings / renames)
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
